$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") for rows 2 through 12: 45208 -> 45212
$ws.Range("C2:C12").Value = 45212
